# Update the quarterly income-statement sheet:
#  - drop the oldest reported quarter (column D) and shift every quarter
#    one column to the left (D<-E, E<-F, ... L<-M)
#  - populate the freed last column (M) with the newly published quarter
#    (also updates the "read_price" derived rows per the new algorithm)
#  - refresh the period-label row (8) and publish-date row (9) headers to
#    match the new rolling 10-quarter window

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- header rows --------------------------------------------------------
# Row 8: quarter labels (D..M), row 9: publish dates (D..M)
$periods = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)
$dates = @(
    "1401-01-15 (3)",
    "1401-03-11 (10)",
    "1401-04-30 (3)",
    "1401-08-18 (4)",
    "1401-10-29 (3)",
    "1402-02-27 (7)",
    "1401-04-30",
    "1401-08-18 (2)",
    "1401-10-29",
    "1402-02-27"
)

$cols = @("D","E","F","G","H","I","J","K","L","M")

$periodsArr = New-Object 'object[,]' 1,10
for ($i = 0; $i -lt 10; $i++) { $periodsArr[0,$i] = $periods[$i] }
$ws.Range("D8:M8").Value = $periodsArr

# Plain "yyyy-mm-dd" strings (no trailing " (n)" suffix) get auto-recognised
# by Excel's Range.Value setter and silently coerced into date serials,
# which would also flip the cell's style. Write those through a
# temporary formula instead and flatten it back to a literal string value
# via copy/paste-special so the cell keeps its original text type & style.
for ($i = 0; $i -lt 10; $i++) {
    $cell = $ws.Range($cols[$i] + "9")
    if ($dates[$i] -match '^\d{4}-\d{2}-\d{2}$') {
        $cell.Formula = '="' + $dates[$i] + '"'
        $cell.Copy()
        $cell.PasteSpecial(-4163) # xlPasteValues
    } else {
        $cell.Value = $dates[$i]
    }
}
$excel.CutCopyMode = 0

# ---- data rows -----------------------------------------------------------
# Each row below is the prior D..M window shifted one quarter to the left,
# with a freshly computed value for the new quarter (M).
$rows = @{
    11 = @(651597, 494945, 957603, 1229447, 1302392, 807869, 1210413, 1449665, 3695117, 3464610)
    12 = @(-369941, -168959, -466387, -562141, -728115, -370997, -602665, -631128, -2355622, -2316989)
    13 = @(281656, 325986, 491216, 667306, 574277, 436872, 607748, 818537, 1339495, 1147621)
    14 = @(-19308, -53829, -38454, -44174, -35158, -60121, -104526, -5399, -61383, -42561)
    15 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    16 = @(9058, -863, 1171, 17015, 11278, 10963, 9, 5658, 18100, 4767)
    17 = @(271406, 271294, 453933, 640147, 550397, 387714, 503231, 818796, 1296212, 1109827)
    18 = @(-39322, -66694, -27405, -67696, -65139, -72152, -25989, -136393, -71576, -103546)
    19 = @(205, -70496, 23, -30937, 30914, -33612, 8510, -377, 713, -11200)
    20 = @(232289, 134104, 426551, 541514, 516172, 281950, 485752, 682026, 1225349, 995081)
    21 = @(-54030, 67176, -96966, -119938, -92238, 60280, -109294, -153435, -275711, 36604)
    22 = @(178259, 201280, 329585, 421576, 423934, 342230, 376458, 528591, 949638, 1031685)
    23 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    24 = @(178259, 201280, 329585, 421576, 423934, 342230, 376458, 528591, 949638, 1031685)
    25 = @(699, 789, 724, 927, 932, 752, 827, 1162, 2087, 1146)
    26 = @(255000, 255000, 455000, 455000, 455000, 455000, 455000, 455000, 455000, 900000)
    27 = @(198, 224, 366, 468, 471, 380, 418, 587, 1055, 1146)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $arr = New-Object 'object[,]' 1,10
    for ($i = 0; $i -lt 10; $i++) { $arr[0,$i] = $vals[$i] }
    $ws.Range("D${r}:M${r}").Value = $arr
}
